# Start edits of Results pages 8-10 on manuscript V1.1
#
# The existing "single vs overlapping QTL" summary table (rows 4-17) keeps
# all of its numbers; only the P31 footnote text in C4 is corrected from
# "5 P31" to "6 P31" (the P-column "single" count used in the write-up).
#
# Below that, a new table (rows 19-37) classifies each ionomics element as
# macronutrient / micronutrient / non-essential analogue / harmful, lists
# its periodic-table group, and compares actual vs. expected single-QTL
# counts (chi-square-style enrichment check), finishing with per-category
# SUM / expected / ratio helper cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- correct the P31 footnote text in the existing table -------------------
$ws.Range("C4").Value = "6 P31 single, 9 overlap"

# --- new section header -----------------------------------------------------
$ws.Range("C19").Value = "single"
$ws.Range("D19").Value = "overlap"

# --- harmful elements --------------------------------------------------------
$ws.Range("A20").Value = "Na"
$ws.Range("B20").Value = "11A"
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = "harmful"

$ws.Range("A21").Value = "Al"
$ws.Range("B21").Value = "3A"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = "harmful"
$ws.Range("I21").Formula = "=77/18"
$ws.Range("J21").Value = "expect 4 QTL per element, on average."

$ws.Range("A22").Value = "Cd"
$ws.Range("B22").Value = "2B"
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = "harmful"
$ws.Range("G22").Value = "actual"
$ws.Range("H22").Value = "expected"

$ws.Range("A23").Value = "As"
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = "harmful"
$ws.Range("G23").Formula = "=SUM(C20:D23)"
$ws.Range("H23").Formula = "=77/(18/4)"
$ws.Range("I23").Formula = "=G23/H23"

# --- macronutrients -----------------------------------------------------------
$ws.Range("A24").Value = "P"
$ws.Range("B24").Value = "5A"
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = "macronutrient"

$ws.Range("A25").Value = "Mg"
$ws.Range("B25").Value = "2A"
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = "macronutrient"

$ws.Range("A26").Value = "K"
$ws.Range("B26").Value = "11A"
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = "macronutrient"

$ws.Range("A27").Value = "Ca"
$ws.Range("B27").Value = "2A"
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = "macronutrient"
$ws.Range("G27").Formula = "=SUM(C24:D27)"
$ws.Range("H27").Formula = "=77/(18/4)"
$ws.Range("I27").Formula = "=G27/H27"

# --- micronutrients -------------------------------------------------------------
$ws.Range("A28").Value = "Cu"
$ws.Range("B28").Value = "1B"
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = "micronutrient"

$ws.Range("A29").Value = "Mn"
$ws.Range("B29").Value = "7B"
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = "micronutrient"

$ws.Range("A30").Value = "Zn"
$ws.Range("B30").Value = "2B"
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = "micronutrient"

$ws.Range("A31").Value = "Fe"
$ws.Range("B31").Value = 8
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = "micronutrient"

$ws.Range("A32").Value = "Mo"
$ws.Range("B32").Value = "6B"
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = "micronutrient"

$ws.Range("A33").Value = "B"
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = "micronutrient"

$ws.Range("A34").Value = "Co"
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = "micronutrient"

$ws.Range("A35").Value = "Se"
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = "micronutrient"
$ws.Range("G35").Formula = "=SUM(C28:D35)"
$ws.Range("H35").Formula = "=77/(18/8)"
$ws.Range("I35").Formula = "=G35/H35"

# --- non-essential analogues ------------------------------------------------------
$ws.Range("A36").Value = "Sr"
$ws.Range("B36").Value = "2A"
$ws.Range("C36").Value = 2
$ws.Range("D36").Value = 7
$ws.Range("E36").Value = "non-essential analogue"

$ws.Range("A37").Value = "Rb"
$ws.Range("B37").Value = "11A"
$ws.Range("C37").Value = 3
$ws.Range("D37").Value = 5
$ws.Range("E37").Value = "non-essential analogue"
$ws.Range("G37").Formula = "=SUM(C36:D37)"
$ws.Range("H37").Formula = "=77/(18/2)"
$ws.Range("I37").Formula = "=G37/H37"

# --- view state: scrolled down to the new table, selection on E32 ----------
$ws.Range("E32").Select()
